# Shiny main page BETA1.0
$wb = $excel.ActiveWorkbook

# --- Data sheet edits ---
$data = $wb.Worksheets.Item("Data")

# B2 was numeric 1 -> becomes text "M"
$data.Range("B2").Value = "M"

# B3 was empty -> becomes text "F"
$data.Range("B3").Value = "F"

# Move the active selection on the Data sheet to B13
$data.Activate()
$null = $data.Range("B13").Select()

$null = $wb.Save()
